$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New week column (L): "SE-45" / "07oct-13nov" -----------------------
# Values first (this extends the used range to column L).
$ws.Range("L2").Value = "SE-45"
$ws.Range("L3").Value = "07oct-13nov"

# Match formatting to the rest of the table by copying each row's existing
# look (column K / B / K) onto the new column L cell(s).
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)

$ws.Range("K2").Copy()
$ws.Range("L2").PasteSpecial(-4122)

$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("L4:L15").PasteSpecial(-4122)

$ws.Range("K16").Copy()
$ws.Range("L16").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Column width / visibility tweak (B narrower, C newly hidden) -------
$ws.Range("C1").EntireColumn.Hidden = $true
$ws.Range("B1").ColumnWidth = 9.3

# --- Selection marker left where the author's cursor ended up -----------
[void]$ws.Range("N14").Select()
